# Replace every "00:00:00" value found in column F (across all worksheets)
# with "Descalificado", matching the race-result "disqualified" status.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # Column F = 6
        if ($cell.Value2 -eq "00:00:00") {
            $cell.Value = "Descalificado"
        }
    }
}
